$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Band1")

# Rows 6-13 correspond to shared strings 63-70 ("Two" through "Nine").
# Add a TRUE boolean value in column B for each of these rows, matching
# the style already used by column A on that row (so the style isn't
# duplicated in the workbook's style table).
for ($r = 6; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    $ws.Cells.Item($r, 2).Value = $true
}
